# "fixed exp month issue and total cost of stay"
# Adds a new guest booking record to the "Victims" sheet (row 2) and
# records the matching Check In / Check Out dates for Room 200 on the
# "Rooms" sheet (row 6).

$wb = $excel.ActiveWorkbook

$victims = $wb.Worksheets.Item("Victims")
$rooms   = $wb.Worksheets.Item("Rooms")

# --- Victims!A2:M2 -------------------------------------------------
# Plain text fields (never look numeric, so Excel keeps them as text
# automatically).
$victims.Range("A2").Value = "John"
$victims.Range("B2").Value = "Doe"
$victims.Range("C2").Value = "jdoe@gmail.com"

# Fields whose text content is purely numeric need an apostrophe
# (quote) prefix so Excel stores them as text instead of silently
# converting them to numbers.
$victims.Range("D2").Value = "'3107957720"
$victims.Range("G2").Value = "'4111111111111111"
$victims.Range("K2").Value = "'200"
$victims.Range("L2").Value = "DOUBLE"
$victims.Range("M2").Value = "'190"

# --- Rooms!B6:C6 (Room 200 stay) -----------------------------------
# Stored as plain text dates (not true date serials), so use the
# quote-prefix trick here too.
$rooms.Range("B6").Value = "'2024-04-26"
$rooms.Range("C6").Value = "'2024-04-28"
